$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data: A3 = "None" (text)
$ws.Range("A3").Value = "None"

# Move the selection to the newly added cell, matching the post-edit state
$ws.Range("A3").Select()
